$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume-change snapshot in columns D and E.
# Values with a single decimal point (e.g. "590.51") are prefixed with a
# leading apostrophe so Excel keeps them as text (matching the original
# inlineStr cells) instead of auto-converting them to numbers; values that
# already contain multiple dots (e.g. "63.154.88") are unambiguous text and
# need no such hint.
$ws.Range("D2").Value = "63.154.88"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.574.20"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'590.51"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").Value = "'144.43"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "'0.106"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "'27.18"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "3.037.18"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "63.038.23"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "2.570.25"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'5.74"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("D24").Value = "'67.77"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  +6.52%  "
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "'469.11"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").Value = "'176.50"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "'18.84"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'1.70"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").Value = "'40.06"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "'157.83"
$ws.Range("E43").Value = "  +4.24%  "
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").Value = "'21.37"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("D46").Value = "'0.633"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D47").Value = "'0.0539"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  +0.00%  "
